$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")

# Row 4 currently holds "test" in the Title column; correct it to "Mr."
$ws.Range("A4").Value = "Mr."

# Move the active selection to A4 to match the saved view state
$ws.Activate()
$ws.Range("A4").Select()
